$wb = $excel.ActiveWorkbook

# Rename the second sheet from "crops_fao" to "fao"
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "fao"

# Fix header: A1 on the "fao" sheet should read "name" (to match crops sheet) instead of "crop"
$ws2.Range("A1").Value = "name"

# Update the selected cell on the active sheet to A3
$ws2.Range("A3").Select()

$wb.Save()
